$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("M2").Value = 2.507757
$ws.Range("N2").Value = 7.523270999999999
$ws.Range("O2").Value = 0.07648041298707947
$ws.Range("P2").Value = 0.07648041298707947
$ws.Range("Q2").Value = 1.833191264975
$ws.Range("R2").Value = 16.498721384775
$ws.Range("S2").Value = 0.001280076748765122
$ws.Range("T2").Value = 0.001280076748765122
$ws.Range("G3").Value = 0.7310083333333334
$ws.Range("H3").Value = 2.193025
$ws.Range("I3").Value = 0.01673731480740535
$ws.Range("J3").Value = 0.01673731480740535
$ws.Range("O3").Value = 0.6219651214303167
$ws.Range("P3").Value = 0.6219651214303167
$ws.Range("Q3").Value = 14.90814423188056
$ws.Range("R3").Value = 134.173298086925
$ws.Range("S3").Value = 0.01041002603660531
$ws.Range("T3").Value = 0.01041002603660531
$ws.Range("G4").Value = 0.7310083333333334
$ws.Range("H4").Value = 2.193025
$ws.Range("I4").Value = 0.01673731480740535
$ws.Range("J4").Value = 0.01673731480740535
$ws.Range("M4").Value = 9.887829999999999
$ws.Range("N4").Value = 29.66349
$ws.Range("O4").Value = 0.3015544655826039
$ws.Range("P4").Value = 0.301554465582604
$ws.Range("Q4").Value = 7.228086128583334
$ws.Range("R4").Value = 65.05277515725
$ws.Range("S4").Value = 0.005047212022034923
$ws.Range("T4").Value = 0.005047212022034925
$ws.Range("I5").Value = 0.8536212576586365
$ws.Range("J5").Value = 0.8536212576586365
$ws.Range("M5").Value = 2.507757
$ws.Range("N5").Value = 7.523270999999999
$ws.Range("O5").Value = 0.07648041298707947
$ws.Range("P5").Value = 0.07648041298707947
$ws.Range("Q5").Value = 93.49474818053999
$ws.Range("R5").Value = 841.4527336248598
$ws.Range("S5").Value = 0.06528530632028269
$ws.Range("T5").Value = 0.06528530632028269
$ws.Range("I6").Value = 0.8536212576586365
$ws.Range("J6").Value = 0.8536212576586365
$ws.Range("O6").Value = 0.6219651214303167
$ws.Range("P6").Value = 0.6219651214303167
$ws.Range("Q6").Value = 760.3315690127132
$ws.Range("R6").Value = 6842.984121114419
$ws.Range("S6").Value = 0.5309226491751535
$ws.Range("T6").Value = 0.5309226491751535
$ws.Range("I7").Value = 0.8536212576586365
$ws.Range("J7").Value = 0.8536212576586365
$ws.Range("M7").Value = 9.887829999999999
$ws.Range("N7").Value = 29.66349
$ws.Range("O7").Value = 0.3015544655826039
$ws.Range("P7").Value = 0.301554465582604
$ws.Range("Q7").Value = 368.6402533826
$ws.Range("R7").Value = 3317.7622804434
$ws.Range("S7").Value = 0.2574133021632004
$ws.Range("T7").Value = 0.2574133021632004
$ws.Range("G8").Value = 5.662136666666666
$ws.Range("H8").Value = 16.98641
$ws.Range("I8").Value = 0.129641427533958
$ws.Range("J8").Value = 0.129641427533958
$ws.Range("M8").Value = 2.507757
$ws.Range("N8").Value = 7.523270999999999
$ws.Range("O8").Value = 0.07648041298707947
$ws.Range("P8").Value = 0.07648041298707947
$ws.Range("Q8").Value = 14.19926286079
$ws.Range("R8").Value = 127.79336574711
$ws.Range("S8").Value = 0.009915029918031644
$ws.Range("T8").Value = 0.009915029918031644
$ws.Range("G9").Value = 5.662136666666666
$ws.Range("H9").Value = 16.98641
$ws.Range("I9").Value = 0.129641427533958
$ws.Range("J9").Value = 0.129641427533958
$ws.Range("O9").Value = 0.6219651214303167
$ws.Range("P9").Value = 0.6219651214303167
$ws.Range("Q9").Value = 115.4733075372411
$ws.Range("R9").Value = 1039.25976783517
$ws.Range("S9").Value = 0.08063244621855781
$ws.Range("T9").Value = 0.08063244621855781
$ws.Range("G10").Value = 5.662136666666666
$ws.Range("H10").Value = 16.98641
$ws.Range("I10").Value = 0.129641427533958
$ws.Range("J10").Value = 0.129641427533958
$ws.Range("M10").Value = 9.887829999999999
$ws.Range("N10").Value = 29.66349
$ws.Range("O10").Value = 0.3015544655826039
$ws.Range("P10").Value = 0.301554465582604
$ws.Range("Q10").Value = 55.98624479676666
$ws.Range("R10").Value = 503.8762031709
$ws.Range("S10").Value = 0.03909395139736858
$ws.Range("T10").Value = 0.03909395139736859
